$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "15% S/LFM+CDN/H:1
25% S+SL/LFM+CDN/H:1
10% S/LFBR+CDN/H:1
20% CR/LFM+CDN/H:2
5% CR/LFINF+CDN/H:1
20% CR+PC/LFM+CDN/H:1
5% MUR/LWAL+CDN/H:1"

$ws.Range("B2").WrapText = $true
$ws.Range("B2").RowHeight = 365

$ws.Range("B2:B11").Select()
